$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '66.939.56'
$ws.Range("E2").Value = '  -1.41%  '
$ws.Range("D3").Value = '3.514.47'
$ws.Range("E3").Value = '  +0.36%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '609.08'
$ws.Range("E5").Value = '  +0.44%  '
$ws.Range("D6").Value = '148.05'
$ws.Range("E6").Value = '  -1.82%  '
$ws.Range("D7").Value = '3.512.92'
$ws.Range("E7").Value = '  +0.36%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '0.479'
$ws.Range("E9").Value = '  -1.64%  '
$ws.Range("D10").Value = '0.143'
$ws.Range("E10").Value = '  -0.84%  '
$ws.Range("D11").Value = '8.06'
$ws.Range("E11").Value = '  +6.61%  '
$ws.Range("D12").Value = '0.423'
$ws.Range("E12").Value = '  -1.69%  '
$ws.Range("D13").Value = '0.0000218'
$ws.Range("E13").Value = '  +1.45%  '
$ws.Range("D14").Value = '32.01'
$ws.Range("E14").Value = '  -0.07%  '
$ws.Range("D15").Value = '4.105.25'
$ws.Range("E15").Value = '  +0.29%  '
$ws.Range("D16").Value = '3.521.95'
$ws.Range("E16").Value = '  +0.80%  '
$ws.Range("D17").Value = '67.016.33'
$ws.Range("E17").Value = '  -1.28%  '
$ws.Range("D18").Value = '0.116'
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("D19").Value = '10.79'
$ws.Range("E19").Value = '  +8.60%  '
$ws.Range("D20").Value = '6.47'
$ws.Range("E20").Value = '  -0.35%  '
$ws.Range("D21").Value = '15.39'
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("D22").Value = '438.30'
$ws.Range("E22").Value = '  -1.77%  '
$ws.Range("D23").Value = '0.610'
$ws.Range("E23").Value = '  -2.27%  '
$ws.Range("D24").Value = '79.74'
$ws.Range("E24").Value = '  +0.82%  '
$ws.Range("D25").Value = '3.647.63'
$ws.Range("E25").Value = '  +0.16%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").Value = '0.0000121'
$ws.Range("E27").Value = '  -3.83%  '
$ws.Range("D28").Value = '9.80'
$ws.Range("E28").Value = '  -1.66%  '
$ws.Range("D29").Value = '8.26'
$ws.Range("E29").Value = '  -3.96%  '
$ws.Range("D30").Value = '2.52'
$ws.Range("E30").Value = '  +0.65%  '
$ws.Range("D31").Value = '1.61'
$ws.Range("E31").Value = '  -2.11%  '
$ws.Range("D32").Value = '0.168'
$ws.Range("E32").Value = '  -1.82%  '
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("D34").Value = '25.62'
$ws.Range("E34").Value = '  +0.09%  '
$ws.Range("D35").Value = '5.98'
$ws.Range("E35").Value = '  -2.60%  '
$ws.Range("D36").Value = '1.81'
$ws.Range("E36").Value = '  -1.93%  '
$ws.Range("D37").Value = '8.09'
$ws.Range("E37").Value = '  +1.31%  '
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.02%  '
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.07%  '
$ws.Range("D40").Value = '175.59'
$ws.Range("E40").Value = '  -0.69%  '
$ws.Range("D41").Value = '0.0896'
$ws.Range("E41").Value = '  -0.18%  '
$ws.Range("D42").Value = '5.41'
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("D43").Value = '2.06'
$ws.Range("E43").Value = '  -11.81%  '
$ws.Range("D44").Value = '0.896'
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("D45").Value = '46.15'
$ws.Range("E45").Value = '  -1.08%  '
$ws.Range("D46").Value = '28.14'
$ws.Range("E46").Value = '  -7.28%  '
$ws.Range("D47").Value = '1.26'
$ws.Range("E47").Value = '  -2.03%  '
$ws.Range("D48").Value = '7.48'
$ws.Range("E48").Value = '  -1.63%  '
$ws.Range("D49").Value = '2.46'
$ws.Range("E49").Value = '  -2.67%  '
$ws.Range("D50").Value = '0.997'
$ws.Range("E50").Value = '  +0.70%  '
$ws.Range("D51").Value = '0.248'
$ws.Range("E51").Value = '  -1.42%  '
